# Refresh crypto Price (col D) and Volume(1h) (col E) values pulled from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.364.48"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.021.01"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.598"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "2.317.77"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.758"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.042.56"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "36.545.82"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.52%  "
$ws.Range("D22").Value = "0.0₃0794"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.29%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -4.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.84%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "1.469.60"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0925"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +33.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "

Write-Output "Updated 88 cells"
